# TC01_Trials_Filter_PubmedID-315.xlsx — "ctdc pubmed id and trial arm"
#
# Adds a new leading "TabName" column (A) that labels each data row as
# CasesTab / FilesTab, replaces the old single Cases query + Stat query
# pair with refreshed Cypher (new property-graph traversal via of_arm/
# of_trial edges) and adds a brand-new Files query row, each with its own
# Stat query copy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Cypher / text payloads (single-quoted here-strings -> no interpolation,
# so backticks and $-signs inside the queries are taken verbatim).
# ---------------------------------------------------------------------

$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
    WHERE a.pubmed_id IN ['31504139'] 
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.pubmed_id IN ['31504139']
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
  WHERE a.pubmed_id IN ['31504139']
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# ---------------------------------------------------------------------
# Row 1 — header. A new "TabName" header is inserted in column A and the
# rest of the former A:D headers shift right into B:E.
# ---------------------------------------------------------------------

$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# ---------------------------------------------------------------------
# Row 2 — Cases tab.
# ---------------------------------------------------------------------

$ws.Range("A2").Value = "CasesTab"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = $casesQuery
$ws.Range("B2").WrapText = $true

$ws.Range("C2").Value = $statQuery
$ws.Range("C2").WrapText = $true

$ws.Range("D2").Value = "TC01_Trials_Filter_PubmedID-315_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC01_Trials_Filter_PubmedID-315_WebData.xlsx"

$ws.Rows.Item(2).RowHeight = 195

# ---------------------------------------------------------------------
# Row 3 — new Files tab.
# ---------------------------------------------------------------------

$ws.Range("A3").Value = "FilesTab"

$ws.Range("B3").Value = $filesQuery
$ws.Range("B3").WrapText = $true

$ws.Range("C3").Value = $statQuery
$ws.Range("C3").WrapText = $true

$ws.Range("D3").Value = "TC01_Trials_Filter_PubmedID-315_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Trials_Filter_PubmedID-315_WebData.xlsx"

$ws.Rows.Item(3).RowHeight = 409.5

# ---------------------------------------------------------------------
# Column widths — first column is new/narrow (auto-fit to the TabName/
# CasesTab/FilesTab labels); B/C share the old wide query width; D/E
# keep the old filename / WebExcel widths.
# ---------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 8.855
$ws.Columns.Item(2).ColumnWidth = 75.855
$ws.Columns.Item(3).ColumnWidth = 75.855
$ws.Columns.Item(4).ColumnWidth = 70.285
$ws.Columns.Item(5).ColumnWidth = 28.57

# ---------------------------------------------------------------------
# View state — scrolled down so row 3 is at the top, with C3 selected
# (matches the author's saved cursor position).
# ---------------------------------------------------------------------

$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
